$d = $word.ActiveDocument

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Worm</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>high</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> health</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Base attack: red</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Middle range</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Low damage</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Power up: shoot a projectile that costs one life and, if it hits, it adds the enemy as a companion (max 3)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">Gameplay: needs to build up a consistent army for being able to deal some damage, </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>but has to be careful with his aim to avoid losing too much health</w:t></w:r></w:p><w:p w:rsidR="00AD62E9" w:rsidRDefault="008F7F70"><w:r w:rsidRPr="003026DA"><w:t>https://www.hongkiat.com/blog/famous-malicious-computer-viruses/</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="00B050"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>Trojan</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="00B050"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>Low health</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="00B050"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>Base attack green-violet</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="00B050"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>Long range</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="00B050"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>High damage</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="00B050"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>Power up: disguise as an enemy for a short period, that makes you undetectectable.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="00B050"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>Gameplay: shoots form distance, and has a stealth ability to reposition himself in case of danger.</w:t></w:r></w:p><w:p w:rsidR="00AD62E9" w:rsidRDefault="00AD62E9"/><w:p><w:pPr><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t xml:space="preserve">Backdoor </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>Medium health</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>Base attack: blue?</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>Short range or melee</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>High damage</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>Powerup: can perform a small blink, that can teleport him through walls</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="00B0F0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>Gameplay</w:t></w:r><w:r><w:rPr><w:color w:val="00B0F0"/></w:rPr><w:t>: try to progress as quickly as possible, avoiding bullets with his high mobility and assassinate them closing the gap between</w:t></w:r></w:p><w:p w:rsidR="008F7F70" w:rsidRDefault="008F7F70"/><w:p><w:pPr><w:rPr><w:color w:val="CCCC00"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="CCCC00"/></w:rPr><w:t>Fisher</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="CCCC00"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="CCCC00"/></w:rPr><w:t>Medium health</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="CCCC00"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="CCCC00"/></w:rPr><w:t>Base attack: yellow</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="CCCC00"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="CCCC00"/></w:rPr><w:t>Middle range</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="CCCC00"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="CCCC00"/></w:rPr><w:t>Medium damage</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:color w:val="CCCC00"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="CCCC00"/></w:rPr><w:lastRenderedPageBreak/><w:t>Powerup: can place baits that will attract enemies</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="CCCC00"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="CCCC00"/></w:rPr><w:t>Gameplay</w:t></w:r><w:r><w:rPr><w:color w:val="CCCC00"/></w:rPr><w:t xml:space="preserve">: proceeds slowly using his powerup for tricky situations. Since is </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="CCCC00"/></w:rPr><w:t>really powerful</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="CCCC00"/></w:rPr><w:t>, it will have a long cd</w:t></w:r></w:p>
'@

$d.Content.InsertXML($xml)
